$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to their rounded (2-decimal) equivalents.
$ws.Range("B5").Value = 15.73
$ws.Range("C5").Value = 11.7
$ws.Range("D5").Value = 1.07
$ws.Range("E5").Value = 34.47
$ws.Range("F5").Value = 27.87
$ws.Range("G5").Value = 12.32
$ws.Range("H5").Value = 48.81
$ws.Range("I5").Value = 19.2
$ws.Range("J5").Value = 8.46
$ws.Range("K5").Value = 12.42
$ws.Range("L5").Value = 13.82
$ws.Range("M5").Value = 14.75
$ws.Range("N5").Value = 3.88
$ws.Range("O5").Value = 12.41
$ws.Range("P5").Value = 17.6
$ws.Range("Q5").Value = 10.58
$ws.Range("R5").Value = 0.74
$ws.Range("S5").Value = 0.69
$ws.Range("T5").Value = 181.52
$ws.Range("U5").Value = 34.71
$ws.Range("V5").Value = 11.45
$ws.Range("W5").Value = 23.2
$ws.Range("X5").Value = 12.16
$ws.Range("Y5").Value = 1.88
$ws.Range("Z5").Value = 23.74
$ws.Range("AA5").Value = 10.12
$ws.Range("AB5").Value = 9.02
$ws.Range("AC5").Value = 10.59
$ws.Range("AD5").Value = 14.53
$ws.Range("AE5").Value = 0.53
$ws.Range("AF5").Value = 44.39
$ws.Range("AG5").Value = 6.38
$ws.Range("AH5").Value = 14.32

# Remove the now-obsolete last row (row 6) entirely.
$ws.Rows.Item(6).Delete()
